# REVER_DailyTracker_NIRMAL.xlsx — add Jan 23-27, 2021 entries to the
# JAN-2021 sheet (rows 24-28), matching the author's latest upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JAN-2021")
$ws.Activate()

# ---------------------------------------------------------------------
# Row 24 - Sat 23 Jan 2021 - Week off
# Built from row 17 (same "Week off" pattern) so B/D pick up the same
# direct formatting (style) that the workbook already uses elsewhere.
# ---------------------------------------------------------------------
$ws.Range("A17").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E23").Copy()
$ws.Range("E24:F24").PasteSpecial(-4122)
$ws.Range("G17").Copy()
$ws.Range("G24").PasteSpecial(-4122)

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 44219
$ws.Range("D24").Value = "Week off"

# ---------------------------------------------------------------------
# Row 25 - Sun 24 Jan 2021 - Week off
# ---------------------------------------------------------------------
$ws.Range("A17").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E23").Copy()
$ws.Range("E25:F25").PasteSpecial(-4122)
$ws.Range("G17").Copy()
$ws.Range("G25").PasteSpecial(-4122)

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 44220
$ws.Range("D25").Value = "Week off"

# ---------------------------------------------------------------------
# Row 26 - Mon 25 Jan 2021 - Task entry (Completed)
# ---------------------------------------------------------------------
$ws.Range("A22:G22").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 44221
$ws.Range("C26").Value = "B2C/B2B app, Sonia & Muji Store"
$ws.Range("D26").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. " + [char]10 + "Regression testing and Retesting on B2C and B2B app." + [char]10 + "Load testing locally on Muji store application"
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = "Completed"
$ws.Rows.Item(26).RowHeight = 60

# ---------------------------------------------------------------------
# Row 27 - Tue 26 Jan 2021 - Holiday
# ---------------------------------------------------------------------
$ws.Range("A16").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E23").Copy()
$ws.Range("E27:F27").PasteSpecial(-4122)
$ws.Range("G16").Copy()
$ws.Range("G27").PasteSpecial(-4122)

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 44222
$ws.Range("D27").Value = "Holiday"

# ---------------------------------------------------------------------
# Row 28 - Wed 27 Jan 2021 - Task entry (Completed)
# ---------------------------------------------------------------------
$ws.Range("A22:G22").Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 44223
$ws.Range("C28").Value = "B2C/B2B app, Sonia & Muji Store"
$ws.Range("D28").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. " + [char]10 + "Regression testing and Retesting on B2C and B2B app." + [char]10 + "Regression testing on Sonia Best ivc report application" + [char]10 + "Load testing locally on Muji store application"
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = "Completed"
$ws.Rows.Item(28).RowHeight = 60

# ---------------------------------------------------------------------
# View state - mirror the author's saved selection/scroll position.
# ---------------------------------------------------------------------
$ws.Range("A33").Select() | Out-Null
